$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 515.56525
$ws.Range("I12").Value = 466.3158
$ws.Range("K12").Value = 466.3158
$ws.Range("M12").Value = -296.3158

$ws.Range("H33").Value = 14860.8
$ws.Range("I33").Value = 17820.25
$ws.Range("J33").Value = 3023
$ws.Range("K33").Value = 17820.25
$ws.Range("L33").Value = 3023
$ws.Range("M33").Value = -17591.25
$ws.Range("N33").Value = -3481

$ws.Range("H54").Value = 14500
$ws.Range("J54").Value = 30000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -30972

$ws.Range("H64").Value = 16600
$ws.Range("I64").Value = 14332.333
$ws.Range("J64").Value = 20001.5
$ws.Range("K64").Value = 14332.333
$ws.Range("L64").Value = 20001.5
$ws.Range("M64").Value = -14084.333
$ws.Range("N64").Value = -20497.5

$ws.Range("H67").Value = 16600
$ws.Range("I67").Value = 14332.333
$ws.Range("J67").Value = 20001.5
$ws.Range("K67").Value = 14332.333
$ws.Range("L67").Value = 20001.5
$ws.Range("M67").Value = -13474.333
$ws.Range("N67").Value = -21717.5

$ws.Range("H98").Value = 729.087
$ws.Range("I98").Value = 729.087
$ws.Range("K98").Value = 729.087
$ws.Range("M98").Value = 768.913

$ws.Range("H122").Value = 729.087
$ws.Range("I122").Value = 729.087
$ws.Range("K122").Value = 2187.261
$ws.Range("M122").Value = 262.739

$ws.Range("H125").Value = 4518.5713
$ws.Range("J125").Value = 5533
$ws.Range("L125").Value = 49797
$ws.Range("N125").Value = -54717

$ws.Range("H132").Value = 3043.739
$ws.Range("I132").Value = 2998.0232
$ws.Range("J132").Value = 3699
$ws.Range("K132").Value = 8994.069600000001
$ws.Range("L132").Value = 11097
$ws.Range("M132").Value = -6464.069600000001
$ws.Range("N132").Value = -16157

$ws.Range("H137").Value = 5782.933
$ws.Range("J137").Value = 8143.5713
$ws.Range("L137").Value = 24430.7139
$ws.Range("N137").Value = -29530.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5102.125
$ws.Range("I2").Value = 3712.5
$ws.Range("K2").Value = 3712.5
$ws.Range("M2").Value = -3599.5

$ws.Range("H32").Value = 13164318
$ws.Range("I32").Value = 15388477
$ws.Range("K32").Value = 15388477
$ws.Range("M32").Value = -15388190

$ws.Range("H61").Value = 25002832
$ws.Range("I61").Value = 30304844
$ws.Range("J61").Value = 7637.2856
$ws.Range("K61").Value = 30304844
$ws.Range("L61").Value = 7637.2856
$ws.Range("M61").Value = -30304632
$ws.Range("N61").Value = -8061.2856

$ws.Range("H116").Value = 5102.125
$ws.Range("I116").Value = 3712.5
$ws.Range("K116").Value = 3712.5
$ws.Range("M116").Value = -1418.5

$ws.Range("H132").Value = 27102708
$ws.Range("I132").Value = 12344.1
$ws.Range("K132").Value = 37032.3
$ws.Range("M132").Value = -34502.3

$ws.Range("H136").Value = 25002832
$ws.Range("I136").Value = 30304844
$ws.Range("J136").Value = 7637.2856
$ws.Range("K136").Value = 90914532
$ws.Range("L136").Value = 22911.8568
$ws.Range("M136").Value = -90911982
$ws.Range("N136").Value = -28011.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5102.125
$ws.Range("I3").Value = 3712.5
$ws.Range("K3").Value = 3712.5
$ws.Range("M3").Value = -3598.5

$ws.Range("H20").Value = 1805
$ws.Range("I20").Value = 1796.5
$ws.Range("K20").Value = 1796.5
$ws.Range("M20").Value = -1549.5

$ws.Range("H99").Value = 4630.091
$ws.Range("I99").Value = 3753.1875
$ws.Range("K99").Value = 3753.1875
$ws.Range("M99").Value = -2255.1875

$ws.Range("H134").Value = 2174.0232
$ws.Range("I134").Value = 1829.7028
$ws.Range("J134").Value = 4297.3335
$ws.Range("K134").Value = 5489.1084
$ws.Range("L134").Value = 12892.0005
$ws.Range("M134").Value = -2954.1084
$ws.Range("N134").Value = -17962.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1116.5
$ws.Range("I16").Value = 862.6667
$ws.Range("J16").Value = 1370.3334
$ws.Range("K16").Value = 862.6667
$ws.Range("L16").Value = 1370.3334
$ws.Range("M16").Value = -575.6667
$ws.Range("N16").Value = -1944.3334

$ws.Range("H31").Value = 83338360
$ws.Range("J31").Value = 125006830
$ws.Range("L31").Value = 125006830
$ws.Range("N31").Value = -125007420

$ws.Range("H34").Value = 83338360
$ws.Range("J34").Value = 125006830
$ws.Range("L34").Value = 125006830
$ws.Range("N34").Value = -125007234

$ws.Range("H113").Value = 1116.5
$ws.Range("I113").Value = 862.6667
$ws.Range("J113").Value = 1370.3334
$ws.Range("K113").Value = 862.6667
$ws.Range("L113").Value = 1370.3334
$ws.Range("M113").Value = 1307.3333
$ws.Range("N113").Value = -5710.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 99.46875
$ws.Range("J2").Value = 142.45454
$ws.Range("L2").Value = 854.7272400000001
$ws.Range("N2").Value = -1080.72724

$ws.Range("H23").Value = 771.06665
$ws.Range("I23").Value = 350
$ws.Range("J23").Value = 876.3333
$ws.Range("K23").Value = 1050
$ws.Range("L23").Value = 2628.9999
$ws.Range("M23").Value = -815
$ws.Range("N23").Value = -3098.9999

$ws.Range("H56").Value = 6485.5835
$ws.Range("I56").Value = 6485.5835
$ws.Range("K56").Value = 6485.5835
$ws.Range("M56").Value = -5955.5835

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 25000
$ws.Range("J44").Value = 25000
$ws.Range("L44").Value = 25000
$ws.Range("N44").Value = -26192

$ws.Range("H52").Value = 44503.5
$ws.Range("J52").Value = 44503.5
$ws.Range("L52").Value = 44503.5
$ws.Range("N52").Value = -45021.5

$ws.Range("H132").Value = 3409.7827
$ws.Range("I132").Value = 3464.7727
$ws.Range("K132").Value = 10394.3181
$ws.Range("M132").Value = -7864.3181

$ws.Range("H136").Value = 47265
$ws.Range("J136").Value = 47265
$ws.Range("L136").Value = 141795
$ws.Range("N136").Value = -146895

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5278.7856
$ws.Range("I7").Value = 4836.1816
$ws.Range("J7").Value = 6901.6665
$ws.Range("K7").Value = 4836.1816
$ws.Range("L7").Value = 6901.6665
$ws.Range("M7").Value = -4724.1816
$ws.Range("N7").Value = -7125.6665

$ws.Range("H16").Value = 2147.875
$ws.Range("I16").Value = 1883.2858
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 1883.2858
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -1713.2858
$ws.Range("N16").Value = -4340

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H46").Value = 2156.7083
$ws.Range("I46").Value = 797.5625
$ws.Range("K46").Value = 797.5625
$ws.Range("M46").Value = -609.5625

$ws.Range("H57").Value = 24000
$ws.Range("J57").Value = 24000
$ws.Range("L57").Value = 24000
$ws.Range("N57").Value = -25132

$ws.Range("H126").Value = 5278.7856
$ws.Range("I126").Value = 4836.1816
$ws.Range("J126").Value = 6901.6665
$ws.Range("K126").Value = 14508.5448
$ws.Range("L126").Value = 20704.9995
$ws.Range("M126").Value = -12038.5448
$ws.Range("N126").Value = -25644.9995

$ws.Range("H132").Value = 62506804
$ws.Range("I132").Value = 4721.28
$ws.Range("J132").Value = 285728540
$ws.Range("K132").Value = 14163.84
$ws.Range("L132").Value = 857185620
$ws.Range("M132").Value = -11633.84
$ws.Range("N132").Value = -857190680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2211.25
$ws.Range("I81").Value = 1365
$ws.Range("K81").Value = 2730
$ws.Range("M81").Value = -1669

$ws.Range("H84").Value = 2211.25
$ws.Range("I84").Value = 1365
$ws.Range("K84").Value = 13650
$ws.Range("M84").Value = -8346

$ws.Range("H132").Value = 2338.7166
$ws.Range("I132").Value = 2213.3877
$ws.Range("K132").Value = 6640.163100000001
$ws.Range("M132").Value = -4110.163100000001
